$wb = $excel.ActiveWorkbook

# 1. Set "EP Service Timeout" header on General Settings sheet first so that
#    it becomes shared string index 140 (before the new JSON string), matching
#    insertion order in the target workbook.
$gs = $wb.Worksheets.Item("General Settings")
$gs.Range("M1").Value = "EP Service Timeout"
$gs.Range("M2").Value = "180"

# 2. Duplicate "VerifyInstall" to create the new "ChangeCustomerConfiguration" sheet.
$src = $wb.Worksheets.Item("VerifyInstall")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item(8)
$newSheet.Name = "ChangeCustomerConfiguration"

$newJson = @"
{
	"customerId": 1001,
	"configuration": {
		"centcom_meta": {
			"schema_version": "1.1.1"
		},
		"global_conf": {
			"log_level": "debug"
		},
		"agent": {
			"ds_host": "endpoint-protection-services.local.tw-test.net",
			"ds_port": 443,
			"ds_protocol": "https",
			"check_update_period": 31,
			"report_period": 35,
			"ds_max_off_period": 24,
			"modules": [
				{
					"name": "Windows Log Monitor",
					"binary_name": "WLM.dll",
					"enabled": true
				},
				{
					"name": "Log File Monitor",
					"binary_name": "LFM.dll",
					"enabled": true
				}
			],
			"transport": {
				"transport_type": 2,
				"syslog": {
					"port": 0
				},
				"scp": {
					"host": "siem-ingress.trustwave.com",
					"dest_folder": "/var/siem/data/nep",
					"port": 9022,
					"user": "twsiem",
					"ack": false,
					"max_send_folder_size": 100
				}
			}
		},
		"wlm": {
			"max_monitor_queue_size": 10000,
			"queues_collector_idle_time": 5,
			"monitor_items": [
				{
					"log_name": "Security",
					"enabled": true,
					"advanced_filter": false,
					"filters": []
				},
				{
					"log_name": "System",
					"enabled": true,
					"advanced_filter": false,
					"filters": []
				}
			]
		},
		"lfm": {
			"max_monitor_queue_size": 10000,
			"queues_collector_idle_time": 5,
			"monitor_items": []
		}
	}
}
"@
$newSheet.Range("A2").Value = $newJson
$newSheet.Range("B2").Value = "180"
$newSheet.Activate()
$newSheet.Range("A2").Select()

# 3. Rename VerifyInstall -> VerifyWinInstall and update its timeout value.
$src.Name = "VerifyWinInstall"
$src.Range("B2").Value = "180"
$src.Range("A2").Select()
